# ----------------------------------------------------------------------------
# Adding data observability code
#   - Numerical sheet: refreshed counts + a new "skewness" column (L)
#   - Categorical sheet: refreshed counts
#   - New "Usage" sheet: table/index usage stats
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Drift
$ws2 = $wb.Worksheets.Item(2)   # Numerical
$ws3 = $wb.Worksheets.Item(3)   # Categorical

# =============================================================================
# Numerical sheet: add "skewness" column header + refreshed values
# =============================================================================
$ws2.Range("L1").Value = "skewness"

$ws2.Range("B2").Value = 15624
$ws2.Range("G2").Value = 15463
$ws2.Range("I2").Value = 37126
$ws2.Range("J2").Value = 19434.18
$ws2.Range("K2").Value = 19926.5
$ws2.Range("L2").Value = -0.09
$ws2.Range("E3").Value = 15624
$ws2.Range("E4").Value = 15624
$ws2.Range("B5").Value = 15624
$ws2.Range("G5").Value = 1882
$ws2.Range("J5").Value = 11665.85
$ws2.Range("L5").Value = 7.33
$ws2.Range("B6").Value = 12457
$ws2.Range("E6").Value = 3167
$ws2.Range("F6").Value = 20.27
$ws2.Range("G6").Value = 5161
$ws2.Range("J6").Value = 6686.31
$ws2.Range("K6").Value = 1840.15
$ws2.Range("L6").Value = 7.99
$ws2.Range("B7").Value = 15624
$ws2.Range("G7").Value = 5325
$ws2.Range("I7").Value = 15718
$ws2.Range("J7").Value = 6158.1
$ws2.Range("K7").Value = 5709
$ws2.Range("L7").Value = 0.32
$ws2.Range("B8").Value = 15624
$ws2.Range("G8").Value = 9950
$ws2.Range("I8").Value = 22189
$ws2.Range("J8").Value = 10561.56
$ws2.Range("K8").Value = 10593
$ws2.Range("L8").Value = 0.05
$ws2.Range("B9").Value = 15624
$ws2.Range("G9").Value = 1406
$ws2.Range("I9").Value = 71501
$ws2.Range("J9").Value = 32665.79
$ws2.Range("L9").Value = -0.34
$ws2.Range("B10").Value = 15624
$ws2.Range("G10").Value = 669
$ws2.Range("J10").Value = 939.88
$ws2.Range("L10").Value = 10.2
$ws2.Range("B11").Value = 15624
$ws2.Range("C11").Value = 7033
$ws2.Range("J11").Value = 20.03
$ws2.Range("L11").Value = 24.16
$ws2.Range("B12").Value = 15624
$ws2.Range("C12").Value = 6089
$ws2.Range("G12").Value = 311
$ws2.Range("J12").Value = 32.95
$ws2.Range("L12").Value = 83.40000000000001
$ws2.Range("B13").Value = 15484
$ws2.Range("J13").Value = 3.45
$ws2.Range("L13").Value = 2.44
$ws2.Range("B14").Value = 15624
$ws2.Range("C14").Value = 15081
$ws2.Range("L14").Value = 5.08
$ws2.Range("B15").Value = 15624
$ws2.Range("C15").Value = 14438
$ws2.Range("J15").Value = 0.08
$ws2.Range("L15").Value = 3.2
$ws2.Range("B16").Value = 15624
$ws2.Range("C16").Value = 15274
$ws2.Range("L16").Value = 6.46
$ws2.Range("B17").Value = 15624
$ws2.Range("C17").Value = 11497
$ws2.Range("J17").Value = 0.26
$ws2.Range("L17").Value = 1.07
$ws2.Range("B18").Value = 15293
$ws2.Range("C18").Value = 14390
$ws2.Range("F18").Value = 2.12
$ws2.Range("L18").Value = 3.74
$ws2.Range("B19").Value = 15293
$ws2.Range("C19").Value = 14390
$ws2.Range("F19").Value = 2.12
$ws2.Range("L19").Value = 3.74
$ws2.Range("B20").Value = 15624
$ws2.Range("C20").Value = 14714
$ws2.Range("L20").Value = 3.77

# Apply the existing bold/bordered header style (style index 2) to the new L1 header
$ws2.Range("K1").Copy()
$ws2.Range("L1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("A1").Select() | Out-Null

# =============================================================================
# Categorical sheet: refreshed values
# =============================================================================
$ws3.Range("B2").Value = 15624
$ws3.Range("B3").Value = 15624
$ws3.Range("B4").Value = 15624
$ws3.Range("B5").Value = 15624
$ws3.Range("B6").Value = 10566
$ws3.Range("D6").Value = 5058
$ws3.Range("E6").Value = 32.37
$ws3.Range("B7").Value = 7063
$ws3.Range("D7").Value = 8561
$ws3.Range("E7").Value = 54.79

# =============================================================================
# New "Usage" sheet (added after Categorical)
# =============================================================================
$wsUsage = $wb.Worksheets.Add($null, $ws3)
$wsUsage.Name = "Usage"

$wsUsage.Range("A1").Value = "schemaname"
$wsUsage.Range("B1").Value = "no_of_times_accessed"
$wsUsage.Range("C1").Value = "table_name"
$wsUsage.Range("D1").Value = "indexrelname"
$wsUsage.Range("E1").Value = "tables_usability"
$wsUsage.Range("F1").Value = "index_usability"

$wsUsage.Range("A2").Value = "adaptiveai"
$wsUsage.Range("B2").Value = 34
$wsUsage.Range("C2").Value = "project_dim"
$wsUsage.Range("E2").Value = "Used"
$wsUsage.Range("F2").Value = "Index not used"

# Apply the same bold/bordered header style to the new sheet's header row
$ws2.Range("A1:F1").Copy()
$wsUsage.Range("A1:F1").PasteSpecial(-4122)   # xlPasteFormats
$wsUsage.Range("A1").Select() | Out-Null

# Restore the originally active sheet/selection
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
